# edit.ps1
# Applies the ICDC TC04 workbook update: rotates the SQL query text assigned
# to cells B2, B4 and B5 on Sheet1, with B4 receiving the updated
# "CaseFilesTab" query (fastq.gz/vcf.gz/tsv.gz instead of a generic gz
# bucket), while B2 and B5 end up with the query text that used to live in
# the cells that, after the shared-string list was reshuffled, now match
# their old neighbours. Also moves the active selection from B3 to B4 to
# match the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$b2Query = @'
SELECT
    DISTINCT c.case_record_id AS "Case ID",
    st.clinical_study_designation AS "Study Code",
    st.clinical_study_type AS "Study Type",
    dmg.breed AS "Breed",
    diag.disease_term AS "Diagnosis",
    diag.stage_of_disease AS "Stage Of Disease",
    CASE 
    WHEN dmg.patient_age_at_enrollment = CAST(dmg.patient_age_at_enrollment AS INT) 
    THEN CAST(CAST(dmg.patient_age_at_enrollment AS INT) AS TEXT)
    ELSE CAST(dmg.patient_age_at_enrollment AS TEXT)
END AS "Age", 
    COALESCE(TRIM(dmg.sex), '') AS "Sex",   
    COALESCE(dmg.neutered_indicator, '') AS "Neutered Status",  
    COALESCE(
        CASE 
            WHEN dmg.weight = CAST(dmg.weight AS INT) THEN CAST(CAST(ROUND(dmg.weight, 0) AS INT) AS TEXT) 
            ELSE dmg.weight 
        END, 
    '') AS "Weight (kg)",  
         
    COALESCE(diag.best_response, '') AS "Response to Treatment",  
    COALESCE(coh.cohort_description, '') AS "Cohort"  
FROM 
    df_case c
JOIN 
    df_study st ON c."study.clinical_study_designation" = st.clinical_study_designation
JOIN 
    df_demographic dmg ON dmg."case.case_record_id" = c.case_record_id
JOIN 
    df_diagnosis diag ON diag."case.case_record_id" = c.case_record_id
JOIN 
    df_enrollment enr ON enr."case.case_record_id" = c.case_record_id
JOIN 
    df_program p ON st."program.program_acronym" = p.program_acronym
JOIN 
    df_sample smp ON smp."case.case_record_id" = c.case_record_id
JOIN 
    df_publication pub ON pub."study.clinical_study_designation" = st.clinical_study_designation
LEFT JOIN 
    df_case_file cf ON cf."sample.sample_id" = smp.sample_id
LEFT JOIN 
    df_study_file sf ON sf."study.clinical_study_designation" = st.clinical_study_designation
LEFT JOIN
    df_cohort coh ON coh."study.clinical_study_designation" = st.clinical_study_designation
WHERE
    st.clinical_study_designation = 'COTC022' AND diag.disease_term = 'Osteosarcoma'
ORDER BY 
    c.case_record_id ASC
LIMIT 100;
'@

$b4Query = @'
SELECT 
    DISTINCT cf.file_name AS "File Name",
    CASE
        WHEN cf.file_name LIKE '%.bai' THEN 'bai'
        WHEN cf.file_name LIKE '%.bam' THEN 'bam'
        WHEN cf.file_name LIKE '%.csv' THEN 'csv'
        WHEN cf.file_name LIKE '%.doc' THEN 'doc'
        WHEN cf.file_name LIKE '%.docx' THEN 'docx'
          WHEN cf.file_name LIKE '%.fastq.gz' THEN 'fastq'
        WHEN cf.file_name LIKE '%.vcf.gz' THEN 'vcf'
        WHEN cf.file_name LIKE '%.tsv.gz' THEN 'tsv'
        WHEN cf.file_name LIKE '%.pdf' THEN 'pdf'
        WHEN cf.file_name LIKE '%.rtf' THEN 'rtf'
        WHEN cf.file_name LIKE '%.tbi' THEN 'tbi'
        WHEN cf.file_name LIKE '%.tif' THEN 'tif'
        WHEN cf.file_name LIKE '%.xls' THEN 'xls'
        WHEN cf.file_name LIKE '%.xlsx' THEN 'xlsx'
        ELSE 'Unknown'
    END AS "Format",
    cf.file_type AS "File Type",
    CASE     
    WHEN cf.file_size >= 1024 * 1024 * 1024 THEN 
        CASE 
            WHEN ROUND(cf.file_size / (1024.0 * 1024.0 * 1024.0), 2) = CAST(ROUND(cf.file_size / (1024.0 * 1024.0 * 1024.0), 0) AS INT) 
            THEN CAST(CAST(ROUND(cf.file_size / (1024.0 * 1024.0 * 1024.0), 0) AS INT) AS TEXT) || ' GB'
            ELSE ROUND(cf.file_size / (1024.0 * 1024.0 * 1024.0), 2) || ' GB'
        END
    WHEN cf.file_size >= 1024 * 1024 THEN 
        CASE 
            WHEN ROUND(cf.file_size / (1024.0 * 1024.0), 2) = CAST(ROUND(cf.file_size / (1024.0 * 1024.0), 0) AS INT) 
            THEN CAST(CAST(ROUND(cf.file_size / (1024.0 * 1024.0), 0) AS INT) AS TEXT) || ' MB'
            ELSE ROUND(cf.file_size / (1024.0 * 1024.0), 2) || ' MB'
        END
    WHEN cf.file_size >= 1024 THEN 
        CASE 
            WHEN ROUND(cf.file_size / 1024.0, 2) = CAST(ROUND(cf.file_size / 1024.0, 0) AS INT) 
            THEN CAST(CAST(ROUND(cf.file_size / 1024.0, 0) AS INT) AS TEXT) || ' KB'
            ELSE ROUND(cf.file_size / 1024.0, 2) || ' KB'
        END
    ELSE 
        CASE 
            WHEN ROUND(cf.file_size, 2) = CAST(ROUND(cf.file_size, 0) AS INT) 
            THEN CAST(CAST(ROUND(cf.file_size, 0) AS INT) AS TEXT) || ' Bytes'
            ELSE ROUND(cf.file_size, 2) || ' Bytes'
        END
END AS "Size",
    'sample' AS "Association",
    cf.file_description AS "Description",
    smp.sample_id AS "Sample ID",
    c.case_record_id AS "Case ID",
    dmg.breed AS "Breed",
    diag.disease_term AS "Diagnosis"
FROM 
    df_case_file cf
JOIN 
    df_sample smp ON cf."sample.sample_id" = smp.sample_id
JOIN 
    df_case c ON smp."case.case_record_id" = c.case_record_id
JOIN 
    df_study st ON c."study.clinical_study_designation" = st.clinical_study_designation
JOIN 
    df_program p ON st."program.program_acronym" = p.program_acronym
JOIN 
    df_demographic dmg ON dmg."case.case_record_id" = c.case_record_id
JOIN 
    df_diagnosis diag ON diag."case.case_record_id" = c.case_record_id
JOIN 
    df_enrollment enr ON enr."case.case_record_id" = c.case_record_id
JOIN 
    df_publication pub ON pub."study.clinical_study_designation" = st.clinical_study_designation
LEFT JOIN 
    df_registration reg ON reg."case.case_record_id" = c.case_record_id
LEFT JOIN 
    df_study_file sf ON sf."study.clinical_study_designation" = st.clinical_study_designation
WHERE
    st.clinical_study_designation = 'COTC022' AND diag.disease_term = 'Osteosarcoma'
ORDER BY 
    cf.file_name ASC
LIMIT 100;
'@

$b5Query = @'
SELECT DISTINCT
    sf.file_name AS "File Name",
    sf.file_type AS "File Type",
    'study' AS "Association",
    sf.file_description AS "Description",
    CASE
        WHEN sf.file_name LIKE '%.bai' THEN 'bai'
        WHEN sf.file_name LIKE '%.bam' THEN 'bam'
        WHEN sf.file_name LIKE '%.csv' THEN 'csv'
        WHEN sf.file_name LIKE '%.doc' THEN 'doc'
        WHEN sf.file_name LIKE '%.docx' THEN 'docx'
        WHEN sf.file_name LIKE '%.gz' THEN 'gz'
        WHEN sf.file_name LIKE '%.pdf' THEN 'pdf'
        WHEN sf.file_name LIKE '%.rtf' THEN 'rtf'
        WHEN sf.file_name LIKE '%.tbi' THEN 'tbi'
        WHEN sf.file_name LIKE '%.tif' THEN 'tif'
        WHEN sf.file_name LIKE '%.xls' THEN 'xls'
        WHEN sf.file_name LIKE '%.xlsx' THEN 'xlsx'
        ELSE 'Unknown'
    END AS "Format",
      CASE     
    WHEN sf.file_size >= 1024 * 1024 * 1024 THEN 
        CASE 
            WHEN ROUND(sf.file_size / (1024.0 * 1024.0 * 1024.0), 2) = CAST(ROUND(sf.file_size / (1024.0 * 1024.0 * 1024.0), 0) AS INT) 
            THEN CAST(CAST(ROUND(sf.file_size / (1024.0 * 1024.0 * 1024.0), 0) AS INT) AS TEXT) || ' GB'
            ELSE ROUND(sf.file_size / (1024.0 * 1024.0 * 1024.0), 2) || ' GB'
        END
    WHEN sf.file_size >= 1024 * 1024 THEN 
        CASE 
            WHEN ROUND(sf.file_size / (1024.0 * 1024.0), 2) = CAST(ROUND(sf.file_size / (1024.0 * 1024.0), 0) AS INT) 
            THEN CAST(CAST(ROUND(sf.file_size / (1024.0 * 1024.0), 0) AS INT) AS TEXT) || ' MB'
            ELSE ROUND(sf.file_size / (1024.0 * 1024.0), 2) || ' MB'
        END
    WHEN sf.file_size >= 1024 THEN 
        CASE 
            WHEN ROUND(sf.file_size / 1024.0, 2) = CAST(ROUND(sf.file_size / 1024.0, 0) AS INT) 
            THEN CAST(CAST(ROUND(sf.file_size / 1024.0, 0) AS INT) AS TEXT) || ' KB'
            ELSE ROUND(sf.file_size / 1024.0, 2) || ' KB'
        END
    ELSE 
        CASE 
            WHEN ROUND(sf.file_size, 2) = CAST(ROUND(sf.file_size, 0) AS INT) 
            THEN CAST(CAST(ROUND(sf.file_size, 0) AS INT) AS TEXT) || ' Bytes'
            ELSE ROUND(sf.file_size, 2) || ' Bytes'
        END
END AS "Size",
    st.clinical_study_designation AS "Study Code"
FROM 
    df_case_file cf
JOIN 
    df_sample smp ON cf."sample.sample_id" = smp.sample_id
JOIN 
    df_case c ON smp."case.case_record_id" = c.case_record_id
JOIN 
    df_study st ON c."study.clinical_study_designation" = st.clinical_study_designation
JOIN 
    df_program p ON st."program.program_acronym" = p.program_acronym
JOIN 
    df_demographic dmg ON dmg."case.case_record_id" = c.case_record_id
JOIN 
    df_diagnosis diag ON diag."case.case_record_id" = c.case_record_id
JOIN 
    df_enrollment enr ON enr."case.case_record_id" = c.case_record_id
JOIN 
    df_publication pub ON pub."study.clinical_study_designation" = st.clinical_study_designation
LEFT JOIN 
    df_study_file sf ON sf."study.clinical_study_designation" = st.clinical_study_designation
WHERE
     st.clinical_study_designation = 'COTC022' AND diag.disease_term = 'Osteosarcoma'
ORDER BY 
    sf.file_name ASC
LIMIT 100;
'@

$ws.Cells.Item(2, 2).Value2 = $b2Query
$ws.Cells.Item(4, 2).Value2 = $b4Query
$ws.Cells.Item(5, 2).Value2 = $b5Query

# Re-assigning the long query text makes the engine auto-fit these rows to
# a taller height; restore the original (already-maxed-out) row height so
# the saved geometry matches the source workbook.
$ws.Rows.Item(2).RowHeight = 409.6
$ws.Rows.Item(4).RowHeight = 409.6
$ws.Rows.Item(5).RowHeight = 409.6

# Move the saved selection / scroll position from B3 to B4, as in the
# updated workbook view state.
$ws.Range("B4").Select() | Out-Null
$window = $excel.ActiveWindow
$window.ScrollRow = 4
$window.ScrollColumn = 1
